# feat: add regionTree api and relation functions
#
# Update the "村居统计" sheet's header row to use the same "*名称" wording
# as the "劝导站统计" sheet, then flip which tab/cell is active & selected
# (村居统计 becomes the active tab; 劝导站统计's selection moves to C1).

$wb = $excel.ActiveWorkbook

$checkpointSheet = $wb.Worksheets.Item(1)   # 劝导站统计
$villageSheet    = $wb.Worksheets.Item(2)   # 村居统计

# --- 村居统计: re-word the header row ---------------------------------
$villageSheet.Range("A1").Value = "县区名称"
$villageSheet.Range("B1").Value = "县区名称"
$villageSheet.Range("C1").Value = "县区名称"

# --- selections -------------------------------------------------------
$checkpointSheet.Range("A1:C1").Select()
$villageSheet.Range("G4").Select()

# --- 村居统计 becomes the active/selected tab --------------------------
$villageSheet.Activate()
